$wb = $excel.ActiveWorkbook

# --- Sheet "Generic" (sheet1) ---------------------------------------------
$wsGeneric = $wb.Worksheets.Item("Generic")
$wsGeneric.Range("B8").Value = 2

# --- Sheet "Productdata" (sheet5) ------------------------------------------
$wsProduct = $wb.Worksheets.Item("Productdata")
$wsProduct.Range("C2").Value = 0
$wsProduct.Range("C3").Value = 0
$wsProduct.Range("C4").Value = 0
$wsProduct.Columns.Item(3).ColumnWidth = 22.90625

# --- Sheet "ForecastedAverageDemand" (sheet6) ------------------------------
$wsForecast = $wb.Worksheets.Item("ForecastedAverageDemand")
$wsForecast.Range("B2").Value = 0
$wsForecast.Range("B3").Value = 0

# New rows 7 and 8, matching the formatting of the existing data rows
$wsForecast.Range("A7").Value = 5
$wsForecast.Range("B7").Value = 10
$wsForecast.Range("C7").Value = 0
$wsForecast.Range("D7").Value = 0

$wsForecast.Range("A8").Value = 6
$wsForecast.Range("B8").Value = 10
$wsForecast.Range("C8").Value = 0
$wsForecast.Range("D8").Value = 0

$wsForecast.Range("A6").Copy()
$wsForecast.Range("A7:A8").PasteSpecial(-4122)

# --- Selections / active sheet ---------------------------------------------
# Order matters: the last sheet activated/selected below becomes the
# workbook's active tab, so Productdata (which should end up tabSelected)
# is activated last.
$wsGeneric.Range("D12").Select()
$wsForecast.Range("B12").Select()

$wsProduct.Activate()
$wsProduct.Range("D10").Select()
